$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values
$ws.Range("A3").Value = -0.876623476148397
$ws.Range("B3").Value = -1.154621533006265

# Update row 9 values
$ws.Range("A9").Value = 0.2752495366249035
$ws.Range("B9").Value = -0.002748520232964513
